# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the zh-cn and de-de handoff rows.
# - Sets the Priority column to "ht" (handoff type) for the rows that were
#   previously blank.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 11, 12, 13, 14)

foreach ($r in $rows) {
    # "Latest HO Xliff Generate Date" column (G) on Overview, shared with the
    # "Latest Handoff Datetime" column (H) on de-de -- both held the same
    # text "2016-08-29 10:19:34" and move to "2016-08-29 10:19:51".
    $overview.Range("G$r").Value = "2016-08-29 10:19:51"
    $dede.Range("H$r").Value = "2016-08-29 10:19:51"

    # "Latest Handoff Datetime" column (H) on zh-cn held
    # "2016-08-29 10:19:30" and moves to "2016-08-29 10:19:46".
    $zhcn.Range("H$r").Value = "2016-08-29 10:19:46"

    # "Priority" column (E) on both locale sheets was blank and is now "ht".
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
